$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column C ("기사제목") so the new
# "발행일" (publish date) column lands in B, shifting the former B/C
# (article title / link) columns right into C/D.
$ws.Range("C1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "수집일"
$ws.Range("B1").Value = "발행일"
$ws.Range("C1").Value = "기사제목"
$ws.Range("D1").Value = "링크"

# Refreshed news rows: A = collected timestamp, B = published timestamp,
# C = article title, D = article link
$data = @(
    @("2026-01-26 04:12:30", "2026-01-26 10:00", "삼성·LG, 프리미엄 AI 노트북 승자는…신제품 국내 출시", "https://n.news.naver.com/mnews/article/003/0013729749?sid=101"),
    @("2026-01-26 04:12:30", "2026-01-26 09:48", "구본욱 KB손보 사장 `"AI 등 가시적 성과 창출해야`"", "https://n.news.naver.com/mnews/article/001/0015866302?sid=101"),
    @("2026-01-26 04:12:30", "2026-01-26 12:01", "갈팡질팡 정부, 결국 ‘대형 원전 건설’로 유턴…AI 전력난 우려에 백...", "https://n.news.naver.com/mnews/article/081/0003611658?sid=102"),
    @("2026-01-26 04:12:30", "2026-01-26 06:00", "''AI 의료기기' 병원 사용 빨라진다… 최단 80일 만에 진입", "https://n.news.naver.com/mnews/article/277/0005711537?sid=101"),
    @("2026-01-26 04:12:30", "2026-01-26 12:50", "삼성 vs. LG, 인공지능(AI) 노트북 신제품 승자는?", "http://woman.chosun.com/news/articleView.html?idxno=124635"),
    @("2026-01-26 04:12:30", "2026-01-26 09:36", "한국한의약진흥원, 한의약 특화 AI 챗봇 ‘맥챗’ 서비스", "https://n.news.naver.com/mnews/article/022/0004100493?sid=102"),
    @("2026-01-26 04:12:30", "2026-01-26 11:34", "카카오 그룹, 지역 AI 인재에 답을 묻다", "https://n.news.naver.com/mnews/article/018/0006207157?sid=105"),
    @("2026-01-26 04:12:30", "2026-01-26 10:17", "수자원공사 `"오픈AI, AI 물관리 한국과 협력 추진 표명`"", "https://n.news.naver.com/mnews/article/421/0008734018?sid=102"),
    @("2026-01-26 04:12:30", "2026-01-26 12:53", "AI·반도체 전력 수요 급증에...탈원전 폐기", "https://n.news.naver.com/mnews/article/016/0002591034?sid=101"),
    @("2026-01-26 04:12:30", "2026-01-26 10:52", "한국수자원공사, OpenAI와 AI 물관리 모색", "https://n.news.naver.com/mnews/article/079/0004108886?sid=102")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
